# "Análise dos dados de 2020"
# Populate Planilha3 with the DATA column (daily dates, most recent first),
# format it like the date columns on the other two sheets, and leave the
# workbook with Planilha3 as the active / selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Planilha3: header + 48 daily dates (2020-06-01 down to 2020-04-15),
#     then two trailing (still date-formatted) blank rows ---
$ws3.Range("A1").Value = "DATA"

$serial = 43983
for ($r = 2; $r -le 49; $r++) {
    $ws3.Cells.Item($r, 1).Value = $serial
    $serial = $serial - 1
}

# Date format (matches the builtin numFmtId 14 "short date" used elsewhere
# in the workbook, but with no alignment override). Set it once, then copy
# the resulting format down through the two empty trailing rows so every
# cell shares a single style record instead of minting one per cell.
$ws3.Range("A2").NumberFormat = "mm-dd-yy"
$ws3.Range("A2").Copy()
$ws3.Range("A3:A51").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A width, like the date columns on the other sheets.
$ws3.Columns.Item(1).ColumnWidth = 12.6

# Selection left on Planilha3 after the last edit.
$ws3.Range("F44").Select()

# --- View tweaks on the other two sheets (scrolled around while
#     reviewing, selection left on Planilha2) ---
$ws2.Range("C47").Select()
$ws2.Application.ActiveWindow.ScrollRow = 28

# --- Planilha3 becomes the active sheet/tab ---
$ws3.Activate()
